# Auto-generated: updates market price / profit data cells per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 141.04762
$ws.Range("I11").Value = 141.04762
$ws.Range("K11").Value = 141.04762
$ws.Range("M11").Value = -1.047619999999995
$ws.Range("H21").Value = 1113.25
$ws.Range("I21").Value = 1113.25
$ws.Range("K21").Value = 1113.25
$ws.Range("M21").Value = -645.25
$ws.Range("H23").Value = 1113.25
$ws.Range("I23").Value = 1113.25
$ws.Range("K23").Value = 1113.25
$ws.Range("M23").Value = -879.25
$ws.Range("H48").Value = 1500
$ws.Range("I48").Value = 1500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -4208
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 1500
$ws.Range("K56").Value = 4500
$ws.Range("M56").Value = -3966
$ws.Range("H64").Value = 3097.6
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 3122.25
$ws.Range("K64").Value = 2999
$ws.Range("L64").Value = 3122.25
$ws.Range("M64").Value = -2751
$ws.Range("N64").Value = -3618.25
$ws.Range("H67").Value = 3097.6
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 3122.25
$ws.Range("K67").Value = 2999
$ws.Range("L67").Value = 3122.25
$ws.Range("M67").Value = -2141
$ws.Range("N67").Value = -4838.25
$ws.Range("H86").Value = 3536.2
$ws.Range("I86").Value = 2920.25
$ws.Range("K86").Value = 2920.25
$ws.Range("M86").Value = -1797.25
$ws.Range("H89").Value = 3536.2
$ws.Range("I89").Value = 2920.25
$ws.Range("K89").Value = 14601.25
$ws.Range("M89").Value = -8985.25
$ws.Range("H112").Value = 2133.6
$ws.Range("J112").Value = 2448.25
$ws.Range("L112").Value = 7344.75
$ws.Range("N112").Value = -9560.75
$ws.Range("H132").Value = 978.37836
$ws.Range("I132").Value = 734.34283
$ws.Range("K132").Value = 2203.02849
$ws.Range("M132").Value = 326.9715099999999
$ws.Range("H141").Value = 2050.4138
$ws.Range("I141").Value = 1587.9286
$ws.Range("K141").Value = 4763.7858
$ws.Range("M141").Value = 416.2142000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1027.25
$ws.Range("I2").Value = 1036.6666
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 1036.6666
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -923.6666
$ws.Range("N2").Value = -1225
$ws.Range("H5").Value = 145.3
$ws.Range("I5").Value = 150.33333
$ws.Range("K5").Value = 150.33333
$ws.Range("M5").Value = -38.33332999999999
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H116").Value = 1027.25
$ws.Range("I116").Value = 1036.6666
$ws.Range("J116").Value = 999
$ws.Range("K116").Value = 1036.6666
$ws.Range("L116").Value = 999
$ws.Range("M116").Value = 1257.3334
$ws.Range("N116").Value = -5587
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1027.25
$ws.Range("I3").Value = 1036.6666
$ws.Range("J3").Value = 999
$ws.Range("K3").Value = 1036.6666
$ws.Range("L3").Value = 999
$ws.Range("M3").Value = -922.6666
$ws.Range("N3").Value = -1227
$ws.Range("H4").Value = 145.3
$ws.Range("I4").Value = 150.33333
$ws.Range("K4").Value = 150.33333
$ws.Range("M4").Value = -35.33332999999999
$ws.Range("H99").Value = 2499.6316
$ws.Range("I99").Value = 2082.5
$ws.Range("K99").Value = 2082.5
$ws.Range("M99").Value = -584.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 417.42856
$ws.Range("I22").Value = 259.5
$ws.Range("K22").Value = 259.5
$ws.Range("M22").Value = 90.5
$ws.Range("H31").Value = 3809.7273
$ws.Range("I31").Value = 3344.6667
$ws.Range("K31").Value = 3344.6667
$ws.Range("M31").Value = -3049.6667
$ws.Range("H34").Value = 3809.7273
$ws.Range("I34").Value = 3344.6667
$ws.Range("K34").Value = 3344.6667
$ws.Range("M34").Value = -3142.6667
$ws.Range("H58").Value = 2179.484
$ws.Range("I58").Value = 1199.9166
$ws.Range("K58").Value = 1199.9166
$ws.Range("M58").Value = -996.9166
$ws.Range("H99").Value = 13094.125
$ws.Range("I99").Value = 9797.700000000001
$ws.Range("K99").Value = 9797.700000000001
$ws.Range("M99").Value = -8299.700000000001
$ws.Range("H122").Value = 3354.6924
$ws.Range("I122").Value = 3575.2727
$ws.Range("K122").Value = 10725.8181
$ws.Range("M122").Value = -8275.8181
$ws.Range("H126").Value = 13094.125
$ws.Range("I126").Value = 9797.700000000001
$ws.Range("K126").Value = 29393.1
$ws.Range("M126").Value = -26923.1
$ws.Range("H127").Value = 80000
$ws.Range("J127").Value = 80000
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920
$ws.Range("H136").Value = 2179.484
$ws.Range("I136").Value = 1199.9166
$ws.Range("K136").Value = 3599.7498
$ws.Range("M136").Value = -1049.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 413.27274
$ws.Range("I122").Value = 308
$ws.Range("K122").Value = 2772
$ws.Range("M122").Value = -322

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 50666.668
$ws.Range("J63").Value = 50666.668
$ws.Range("L63").Value = 50666.668
$ws.Range("N63").Value = -52038.668
$ws.Range("H66").Value = 50666.668
$ws.Range("J66").Value = 50666.668
$ws.Range("L66").Value = 152000.004
$ws.Range("N66").Value = -158864.004
$ws.Range("H100").Value = 500500
$ws.Range("J100").Value = 500500
$ws.Range("L100").Value = 500500
$ws.Range("N100").Value = -502664
$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 1550
$ws.Range("K113").Value = 1550
$ws.Range("M113").Value = 620
$ws.Range("H122").Value = 41135.58
$ws.Range("I122").Value = 2525.4092
$ws.Range("J122").Value = 253491.5
$ws.Range("K122").Value = 7576.2276
$ws.Range("L122").Value = 760474.5
$ws.Range("M122").Value = -5126.2276
$ws.Range("N122").Value = -765374.5
$ws.Range("H123").Value = 55428.57
$ws.Range("J123").Value = 55428.57
$ws.Range("L123").Value = 55428.57
$ws.Range("N123").Value = -60328.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4890.5
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H9").Value = 68.833336
$ws.Range("I9").Value = 68.833336
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 68.833336
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 155.166664
$ws.Range("H16").Value = 3001
$ws.Range("J16").Value = 3001
$ws.Range("L16").Value = 3001
$ws.Range("N16").Value = -3341
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10450
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11560
$ws.Range("H93").Value = 2080.8
$ws.Range("I93").Value = 1976
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1976
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -728
$ws.Range("N93").Value = -4996
$ws.Range("H126").Value = 4890.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21248
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66240
$ws.Range("H122").Value = 2833.3333
$ws.Range("I122").Value = 2600
$ws.Range("K122").Value = 7800
$ws.Range("M122").Value = -5350

